$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2:D13 currently hold the patch/release date as text (shared strings
# such as "25-12-2018"). Replace them with real Excel date serial values
# and format the column as dates, one calendar day apart starting
# 2018-12-26.
$dates = @(43460,43461,43462,43463,43464,43465,43466,43467,43468,43469,43470,43471)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dates[$i]
}

# Apply the date number format to D2 then fan it out to the rest of the
# column so every cell shares a single style entry instead of minting one
# per cell.
$ws.Cells.Item(2, 4).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2, 4).Copy()
$ws.Range("D3:D13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to C4 (matches the saved sheetView state).
$ws.Range("C4").Select()
